$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Concept-Description of the DB:
# u_surname (D7) / u_forename (D8) column size changes from VARCHAR(45) to VARCHAR(35)
$ws.Range("D7").Value = "VARCHAR(35)"
$ws.Range("D8").Value = "VARCHAR(35)"

# C7:C9 (u_surname, u_forename, u_birtdate) now get the same highlighted
# formatting already used by the rest of column C (e.g. C3:C6, C10:C13) -
# copy that format from C6, a representative already-highlighted cell.
$ws.Range("C6").Copy()
$ws.Range("C7:C9").PasteSpecial(-4122)  # xlPasteFormats

# Selection moves from F8 to C7:C8, with C7 as the active cell
$ws.Range("C7:C8").Select()
